$wb = $excel.ActiveWorkbook

# "Croatia" is currently the last sheet in the workbook; duplicate it
# (mirrors Excel's "Move or Copy... > Create a copy") and place the
# duplicate right after it.
$source = $wb.Worksheets.Item("Croatia")
$source.Copy($null, $source)

# The newly created copy is the sheet right after the source.
$newSheet = $wb.Worksheets.Item($source.Index + 1)
$newSheet.Name = "Greece"

# Update the market name / ticket reference cells for the new sheet.
$newSheet.Range("B2").Value = "Greece Market"
$newSheet.Range("B4").Value = "NGC-4119/T3205"

# Restore the source sheet's selection to a full-sheet selection (as
# happens after Excel performs the copy) and make sure the new sheet
# stays the active tab.
$source.Cells.Select()
$newSheet.Activate()
